$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted above the current row 12, shifting the
# existing rows 12-28 down to 13-29 (dimension grows from A1:R28 to A1:R29).
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new week's data.
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Terminal La Palmera de La Serena"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44533
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 100114007
$ws.Range("G12").Value = "Jengibre"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 520
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17500
$ws.Range("N12").Value = "`$/caja 13 kilos"
$ws.Range("O12").Value = "Perú"
$ws.Range("P12").Value = 1346
$ws.Range("Q12").Value = 13
$ws.Range("R12").Value = "Hortaliza"
